$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-by-cell updates per the diff. D-column values are price strings that
# can look numeric (e.g. '219.10', '4.560'); Excel's normal smart-typing would
# coerce them into actual numbers and silently drop significant trailing
# zeros (219.10 -> 219.1). The source workbook stores every Price/Volume cell
# as literal text, so for D-column writes we force Text format first and
# restore the Normal style afterwards (matches the unstyled original cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.335.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5341"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2662"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06397"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07853"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.560"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.680.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.894.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5539"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8192"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.358.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.686"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.044"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1231"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.212"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  +4.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05867"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.283"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.284"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9707"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5827"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01601"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8627"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.066.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.838"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.011"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.806.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.013"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4392"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.009"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05164"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
